$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# "Back" button renamed to "Cancel" button; the dedicated "view timesheet page"
# no longer exists since timesheet management was consolidated into a single
# webpart, so the wording is generalized to "timesheet page".
# NOTE: execution order below matters - it controls how the new shared
# strings are appended/ordered on save.
$ws.Range("D33").Value = 'User is getting redirected to timesheet listing page by clicking on "Cancel" button.'
$ws.Range("B33").Value = 'User should be redirected to timesheet listing page by clicking on "Cancel" button.'
$ws.Range("C33").Value = 'User should be redirected to timesheet listing page by clicking on "Cancel" button.'
$ws.Range("B32").Value = 'User should get "Cancel" button on timesheet page.'
$ws.Range("C32").Value = 'User should get "Cancel" button on timesheet page.'
$ws.Range("D32").Value = 'User is getting get "Cancel" button on timesheet page.'

# Update the active selection on the sheet to reflect where the edit was made.
$ws.Range("B24").Select()
